$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 26608.074
$ws.Range("I28").Value = 44215.75
$ws.Range("J28").Value = 996.9091
$ws.Range("K28").Value = 44215.75
$ws.Range("L28").Value = 996.9091
$ws.Range("M28").Value = -43730.75
$ws.Range("N28").Value = -1966.9091
$ws.Range("H40").Value = 1940.5555
$ws.Range("I40").Value = 1894.6154
$ws.Range("J40").Value = 2060
$ws.Range("K40").Value = 1894.6154
$ws.Range("L40").Value = 2060
$ws.Range("M40").Value = -1719.6154
$ws.Range("N40").Value = -2410
$ws.Range("H62").Value = 2343.5715
$ws.Range("I62").Value = 2601.25
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 2601.25
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -1977.25
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 2343.5715
$ws.Range("I65").Value = 2601.25
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 13006.25
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -9886.25
$ws.Range("N65").Value = -16240
$ws.Range("H110").Value = 48000
$ws.Range("J110").Value = 48000
$ws.Range("L110").Value = 48000
$ws.Range("N110").Value = -56180
$ws.Range("H138").Value = 2230.7742
$ws.Range("I138").Value = 1296.6666
$ws.Range("J138").Value = 5433.4287
$ws.Range("K138").Value = 3889.9998
$ws.Range("L138").Value = 16300.2861
$ws.Range("M138").Value = 1250.0002
$ws.Range("N138").Value = -26580.2861

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 29252.375
$ws.Range("J23").Value = 15667.833
$ws.Range("L23").Value = 15667.833
$ws.Range("N23").Value = -16185.833

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 1200
$ws.Range("I8").Value = 1200
$ws.Range("K8").Value = 1200
$ws.Range("M8").Value = -1060
$ws.Range("H35").Value = 21200
$ws.Range("I35").Value = 6000
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 6000
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -5690
$ws.Range("N35").Value = -25620

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 46500
$ws.Range("J28").Value = 46500
$ws.Range("L28").Value = 46500
$ws.Range("N28").Value = -46990
$ws.Range("H31").Value = 20004920
$ws.Range("I31").Value = 76924880
$ws.Range("J31").Value = 6014
$ws.Range("K31").Value = 76924880
$ws.Range("L31").Value = 6014
$ws.Range("M31").Value = -76924585
$ws.Range("N31").Value = -6604
$ws.Range("H34").Value = 20004920
$ws.Range("I34").Value = 76924880
$ws.Range("J34").Value = 6014
$ws.Range("K34").Value = 76924880
$ws.Range("L34").Value = 6014
$ws.Range("M34").Value = -76924678
$ws.Range("N34").Value = -6418
$ws.Range("H86").Value = 2253.7856
$ws.Range("I86").Value = 1817.125
$ws.Range("J86").Value = 2836
$ws.Range("K86").Value = 1817.125
$ws.Range("L86").Value = 2836
$ws.Range("M86").Value = -694.125
$ws.Range("N86").Value = -5082
$ws.Range("H89").Value = 2253.7856
$ws.Range("I89").Value = 1817.125
$ws.Range("J89").Value = 2836
$ws.Range("K89").Value = 9085.625
$ws.Range("L89").Value = 14180
$ws.Range("M89").Value = -3469.625
$ws.Range("N89").Value = -25412
$ws.Range("H132").Value = 55564772
$ws.Range("I132").Value = 66676390
$ws.Range("J132").Value = 6666
$ws.Range("K132").Value = 200029170
$ws.Range("L132").Value = 19998
$ws.Range("M132").Value = -200026640
$ws.Range("N132").Value = -25058

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 360.72
$ws.Range("I5").Value = 296.09525
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 888.28575
$ws.Range("L5").Value = 2100
$ws.Range("M5").Value = -776.28575
$ws.Range("N5").Value = -2324
$ws.Range("H68").Value = 861.7683
$ws.Range("I68").Value = 533.3090999999999
$ws.Range("J68").Value = 1530.8518
$ws.Range("K68").Value = 1599.9273
$ws.Range("L68").Value = 4592.555399999999
$ws.Range("M68").Value = -788.9272999999998
$ws.Range("N68").Value = -6214.555399999999
$ws.Range("H71").Value = 861.7683
$ws.Range("I71").Value = 533.3090999999999
$ws.Range("J71").Value = 1530.8518
$ws.Range("K71").Value = 4799.7819
$ws.Range("L71").Value = 13777.6662
$ws.Range("M71").Value = -743.7819
$ws.Range("N71").Value = -21889.6662
$ws.Range("H75").Value = 111112290
$ws.Range("I75").Value = 125000696
$ws.Range("J75").Value = 5000
$ws.Range("K75").Value = 375002088
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -375001090
$ws.Range("N75").Value = -16996
$ws.Range("H78").Value = 111112290
$ws.Range("I78").Value = 125000696
$ws.Range("J78").Value = 5000
$ws.Range("K78").Value = 1125006264
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = -1125001272
$ws.Range("N78").Value = -54984
$ws.Range("H80").Value = 3305.95
$ws.Range("I80").Value = 2798.4
$ws.Range("J80").Value = 3475.1333
$ws.Range("K80").Value = 8395.200000000001
$ws.Range("L80").Value = 10425.3999
$ws.Range("M80").Value = -7459.200000000001
$ws.Range("N80").Value = -12297.3999
$ws.Range("H83").Value = 3305.95
$ws.Range("I83").Value = 2798.4
$ws.Range("J83").Value = 3475.1333
$ws.Range("K83").Value = 25185.6
$ws.Range("L83").Value = 31276.1997
$ws.Range("M83").Value = -20505.6
$ws.Range("N83").Value = -40636.1997
$ws.Range("H107").Value = 341.87234
$ws.Range("I107").Value = 199.25581
$ws.Range("J107").Value = 1875
$ws.Range("K107").Value = 597.76743
$ws.Range("L107").Value = 5625
$ws.Range("M107").Value = 1322.23257
$ws.Range("N107").Value = -9465
$ws.Range("H135").Value = 360.72
$ws.Range("I135").Value = 296.09525
$ws.Range("J135").Value = 700
$ws.Range("K135").Value = 2664.85725
$ws.Range("L135").Value = 6300
$ws.Range("M135").Value = -129.85725
$ws.Range("N135").Value = -11370

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 233
$ws.Range("I13").Value = 233
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 233
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -94
$ws.Range("N13").ClearContents()
$ws.Range("H104").Value = 48000
$ws.Range("J104").Value = 48000
$ws.Range("L104").Value = 48000
$ws.Range("N104").Value = -54988
$ws.Range("H107").Value = 296.66666
$ws.Range("I107").Value = 296.4
$ws.Range("K107").Value = 296.4
$ws.Range("M107").Value = 1623.6
$ws.Range("H113").Value = 4319.769
$ws.Range("I113").Value = 3922.182
$ws.Range("J113").Value = 6506.5
$ws.Range("K113").Value = 3922.182
$ws.Range("L113").Value = 6506.5
$ws.Range("M113").Value = -1752.182
$ws.Range("N113").Value = -10846.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1320.1
$ws.Range("I46").Value = 825.25
$ws.Range("J46").Value = 1650
$ws.Range("K46").Value = 825.25
$ws.Range("L46").Value = 1650
$ws.Range("M46").Value = -637.25
$ws.Range("N46").Value = -2026
$ws.Range("H132").Value = 4749.136
$ws.Range("I132").Value = 5217.4165
$ws.Range("J132").Value = 3500.389
$ws.Range("K132").Value = 15652.2495
$ws.Range("L132").Value = 10501.167
$ws.Range("M132").Value = -13122.2495
$ws.Range("N132").Value = -15561.167

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1401
$ws.Range("I126").Value = 1252
$ws.Range("J126").Value = 1550
$ws.Range("K126").Value = 3756
$ws.Range("L126").Value = 4650
$ws.Range("M126").Value = -1286
$ws.Range("N126").Value = -9590
$ws.Range("H136").Value = 3354.9216
$ws.Range("I136").Value = 4103.6562
$ws.Range("J136").Value = 2093.8948
$ws.Range("K136").Value = 12310.9686
$ws.Range("L136").Value = 6281.6844
$ws.Range("M136").Value = -9760.9686
$ws.Range("N136").Value = -11381.6844
